$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 144, pushing existing rows 144:150 down to 145:151
$ws.Rows.Item(144).Insert()

# Copy the number format (date style) used by column D from the row below (now row 145)
$ws.Range("D144").NumberFormat = $ws.Range("D145").NumberFormat

# Populate the new row 144 with the new record's data
$ws.Cells.Item(144, 1).Value = 4
$ws.Cells.Item(144, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(144, 3).Value = "Los Lagos"
$ws.Cells.Item(144, 4).Value = 44516
$ws.Cells.Item(144, 5).Value = 10
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100109
$ws.Cells.Item(144, 8).Value = "Uva"
$ws.Cells.Item(144, 9).Value = 100109001
$ws.Cells.Item(144, 10).Value = "Uva"
$ws.Cells.Item(144, 11).Value = "Superior Seedless"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 200
$ws.Cells.Item(144, 14).Value = 29000
$ws.Cells.Item(144, 15).Value = 30000
$ws.Cells.Item(144, 16).Value = 29500
$ws.Cells.Item(144, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(144, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(144, 19).Value = 2950
$ws.Cells.Item(144, 20).Value = 10
